# Fixed update to excel issue
#
# 1) Rename the "Requested quantity" header on the "Weekly Quantity" sheet
#    to "Weekly_PO_Qty".
# 2) Rename the "Requested quantity" header on the "Monthly Trend" sheet
#    to "Monthly_PO_Qty".
# 3) Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast
#    data (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1) Weekly Quantity header -------------------------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2) Monthly Trend header ----------------------------------------------
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3) New "PO Forecast" sheet -------------------------------------------
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Reuse the existing header style (bold, bordered, centered) from the
# "Weekly Quantity" sheet so the new header row matches the rest of the
# workbook, and reuse the date-formatted style for the "ds" column.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(45578.99999999999, 120, 119.9999929314775, 119.9999929321358),
    @(45592.99999999999, 24, 23.99998093253115, 23.99998093316039),
    @(45599.99999999999, 0, -24.0000253752479, -24.00002475769405),
    @(45606.99999999999, 0, -72.00003212809777, -72.00003003536231),
    @(45613.99999999999, 0, -120.0000391146809, -120.0000349352976),
    @(45620.99999999999, 0, -168.0000461406963, -168.0000398870307),
    @(45627.99999999999, 0, -216.0000534186464, -216.0000444410176),
    @(45634.99999999999, 0, -264.0000608587602, -264.0000489276317),
    @(45641.99999999999, 0, -312.0000686278593, -312.0000531740021),
    @(45648.99999999999, 0, -360.0000762241912, -360.000057398054)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Leave the original sheet selected / active, as before.
$wsWeekly.Activate()
